$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.93%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.22%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.071"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.30%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07925"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.78%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.115"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.13%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.960"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.48%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.136"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.02%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9238"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.35%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09697"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.47%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1846"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.01%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08590"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.44%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03574"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.51%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09943"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.31%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001430"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-4.47%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005718"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.39%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.462"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.05%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "21.97%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.95%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.44%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.175"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.74%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2210"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.48%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04560"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.89%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001237"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.54%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004886"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "9.45%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.91%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004755"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.08%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01848"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.07%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04733"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.22%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007898"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.81%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1398"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007600"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.74%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002193"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.59%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.50%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006300"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.54%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.25%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.29%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "50.94"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "375.70%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-25.58%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002104"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.25%"
